# Update "想去人数" (interest counts) for both the 展览 sheet and the
# 全部类型 sheet, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 3,4,5,7,9 get updated counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 312
$wsExhibit.Range("F4").Value = 224
$wsExhibit.Range("F5").Value = 2722
$wsExhibit.Range("F7").Value = 369
$wsExhibit.Range("F9").Value = 963

# Sheet "全部类型" - rows 3,4,5,7,10 get updated counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 312
$wsAll.Range("F4").Value = 224
$wsAll.Range("F5").Value = 2722
$wsAll.Range("F7").Value = 369
$wsAll.Range("F10").Value = 963
